$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 65, shifting rows 65:189 down to 66:190.
$ws.Rows.Item(65).Insert()

# Populate the newly inserted row 65 with its data. Most fields mirror the
# record that used to occupy row 65 before the shift (now at row 66); only
# the Fecha (D) and Volumen (J) differ, per the commit's new weekly entry.
$ws.Range("A65").Value = 4
$ws.Range("B65").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C65").Value = "Los Lagos"
$ws.Range("D65").Value = 44581
$ws.Range("E65").Value = 10
$ws.Range("F65").Value = 100112017
$ws.Range("G65").Value = "Apio"
$ws.Range("H65").Value = "Americana (o)"
$ws.Range("I65").Value = "Primera"
$ws.Range("J65").Value = 25
$ws.Range("K65").Value = 12000
$ws.Range("L65").Value = 12000
$ws.Range("M65").Value = 12000
$ws.Range("N65").Value = "$/docena de matas"
$ws.Range("O65").Value = "Región de Coquimbo"
$ws.Range("P65").Value = 2000
$ws.Range("Q65").Value = 6
$ws.Range("R65").Value = "Hortaliza"
